$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the sheet: this becomes a "Lista Oficial" (official roster) print
# layout instead of the old plain table.
$ws.Cells.ClearContents()

# Pre-create the merged ranges before touching value/font so Excel does not
# fan the anchor cell's style out into a batch of brand-new style records
# for every merged member cell.
$ws.Range("D4:I4").Merge()
$ws.Range("D5:E5").Merge()
$ws.Range("D6:I6").Merge()

# --- Title block -------------------------------------------------------
# Row 1: "Lista Oficial" big bold underlined title (E1)
$ws.Range("E1").Value = "Lista Oficial"
$ws.Range("E1").Font.Name = "Calibri"
$ws.Range("E1").Font.Size = 18
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Font.Underline = $true

# Row 2: "SICEP - Centro de Computo" bold subtitle (D2)
$ws.Range("D2").Value = "SICEP - Centro de Computo"
$ws.Range("D2").Font.Name = "Calibri"
$ws.Range("D2").Font.Size = 16
$ws.Range("D2").Font.Bold = $true

# Row 4: PROFESOR label + value
$ws.Range("B4").Value = "PROFESOR:"
$ws.Range("D4").Value = "Christian RCGS DEVELOPER"
$ws.Range("D4").Font.Name = "Calibri"
$ws.Range("D4").Font.Size = 14
$ws.Range("D4").Font.Bold = $true

# Row 5: No TRABAJADOR label + value
$ws.Range("B5").Value = "No TRABAJADOR:"
$ws.Range("D5").Value = 312260633
$ws.Range("D5").Font.Name = "Calibri"
$ws.Range("D5").Font.Size = 12
$ws.Range("D5").Font.Bold = $true

# Row 6: SEMINARIO label + value
$ws.Range("B6").Value = "SEMINARIO"
$ws.Range("D6").Value = "003 Aspel COI"
$ws.Range("D6").Font.Name = "Calibri"
$ws.Range("D6").Font.Size = 14
$ws.Range("D6").Font.Bold = $true

# Row 7: GRUPO / SEMESTRE / GENERACION / MODALIDAD labels + values
$ws.Range("B7").Value = "GRUPO:"
$ws.Range("C7").Value = 1000
$ws.Range("C7").Font.Name = "Calibri"
$ws.Range("C7").Font.Size = 12
$ws.Range("C7").Font.Bold = $true

$ws.Range("D7").Value = "SEMESTRE:"

$ws.Range("E7").Value = "2021-2"
$ws.Range("E7").Font.Name = "Calibri"
$ws.Range("E7").Font.Size = 12
$ws.Range("E7").Font.Bold = $true

$ws.Range("F7").Value = "GENERACION:"

$ws.Range("G7").Value = 2022
$ws.Range("G7").Font.Name = "Calibri"
$ws.Range("G7").Font.Size = 12
$ws.Range("G7").Font.Bold = $true

$ws.Range("H7").Value = "MODALIDAD:"

$ws.Range("I7").Value = "Presencial"
$ws.Range("I7").Font.Name = "Calibri"
$ws.Range("I7").Font.Size = 12
$ws.Range("I7").Font.Bold = $true

# --- Student roster rows ------------------------------------------------
# Row 10: new first student
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 55555555
$ws.Range("C10").Value = "Abasolo"
$ws.Range("D10").Value = "Lopez"
$ws.Range("E10").Value = "Juan"
$ws.Range("F10").Value = "juan@algo.com"
$ws.Range("G10").Value = "'5512457845"
$ws.Range("H10").Value = "Ing. en Alimentos"

# Row 11: previously existing student record, pushed down one row
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 456156165
$ws.Range("C11").Value = "Hernandez"
$ws.Range("D11").Value = "Fernandez"
$ws.Range("E11").Value = "Alberto"
$ws.Range("F11").Value = "cuyo@gmail.com"
$ws.Range("G11").Value = "'45641564165"
$ws.Range("H11").Value = "Informatica"

# Selection / view state mirrors the authored workbook
$ws.Range("I7").Select() | Out-Null

Write-Host "done"
